$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New pedigree data for rows 2-6 (columns A,B,C,D,H,I,K,L,M)
# Columns N and O are cleared for all rows (sire_sw_id / dam_sw_id no longer populated)

$data = @(
    @{Row=2; A="3403"; B="73512_1"; C="ALBINO";     D="M"; H="72905_2"; I="72801_1"; K="WHSM0517B"; L="3000"; M="LBRB"},
    @{Row=3; A="3404"; B="73505_1"; C="BLACK";      D="M"; H="73309_2"; I="73347_2"; K="WHSM0561B"; L="3000"; M="RC"},
    @{Row=4; A="3469"; B="73504_4"; C="BRN HOOD";   D="F"; H="72945_1"; I="72797_1"; K="WHSF0507B"; L="3000"; M="LTRT"},
    @{Row=5; A="3470"; B="73511_1"; C="BLACK HOOD"; D="F"; H="73337_1"; I="73332_1"; K="WHSF0551B"; L="3000"; M="LTRC"},
    @{Row=6; A="3471"; B="73546_1"; C="BRN HOOD";   D="F"; H="72794_2"; I="72790_1"; K="WHSF0516B"; L="3000"; M="LTRB"}
)

foreach ($row in $data) {
    $r = $row.Row
    $ws.Range("A$r").Value = "'" + $row.A
    $ws.Range("A$r").ClearFormats()
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Value = $row.I
    $ws.Range("K$r").Value = $row.K
    $ws.Range("L$r").Value = "'" + $row.L
    $ws.Range("L$r").ClearFormats()
    $ws.Range("M$r").Value = $row.M
    $ws.Range("N$r").ClearContents()
    $ws.Range("O$r").ClearContents()
}
